# Scheduled runner update: refresh cached market price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a handful of leve
# rows across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 1000
$ws.Range("M64").Value = -752

$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 1000
$ws.Range("M67").Value = -142

$ws.Range("H107").Value = 658.1905
$ws.Range("J107").Value = 498.33334
$ws.Range("L107").Value = 498.33334
$ws.Range("N107").Value = -4338.33334

$ws.Range("H137").Value = 1028
$ws.Range("I137").Value = 819.35297
$ws.Range("K137").Value = 2458.05891
$ws.Range("M137").Value = 91.9410899999998

$ws.Range("H138").Value = 3229.125
$ws.Range("J138").Value = 3319.4285
$ws.Range("L138").Value = 9958.2855
$ws.Range("N138").Value = -20238.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2028.75
$ws.Range("I45").Value = 2070.4285
$ws.Range("K45").Value = 2070.4285
$ws.Range("M45").Value = -1693.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 489
$ws.Range("I16").Value = 489
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 489
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -202
$ws.Range("N16").ClearContents()

$ws.Range("H33").Value = 26838.143
$ws.Range("I33").Value = 1967.25
$ws.Range("J33").Value = 59999.332
$ws.Range("K33").Value = 1967.25
$ws.Range("L33").Value = 59999.332
$ws.Range("M33").Value = -1588.25
$ws.Range("N33").Value = -60757.332

$ws.Range("H36").Value = 27154.75
$ws.Range("I36").Value = 11747
$ws.Range("K36").Value = 11747
$ws.Range("M36").Value = -11359

$ws.Range("H40").Value = 27154.75
$ws.Range("I40").Value = 11747
$ws.Range("K40").Value = 11747
$ws.Range("M40").Value = -11587

$ws.Range("H62").Value = 4499.8
$ws.Range("I62").Value = 3333
$ws.Range("K62").Value = 3333
$ws.Range("M62").Value = -2709

$ws.Range("H65").Value = 4499.8
$ws.Range("I65").Value = 3333
$ws.Range("K65").Value = 16665
$ws.Range("M65").Value = -13545

$ws.Range("H74").Value = 86875
$ws.Range("J74").Value = 86875
$ws.Range("L74").Value = 86875
$ws.Range("N74").Value = -88623

$ws.Range("H77").Value = 86875
$ws.Range("J77").Value = 86875
$ws.Range("L77").Value = 260625
$ws.Range("N77").Value = -269361

$ws.Range("H113").Value = 489
$ws.Range("I113").Value = 489
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 489
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1681
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 2155.889
$ws.Range("I132").Value = 2050.375
$ws.Range("K132").Value = 6151.125
$ws.Range("M132").Value = -3621.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 250.38461
$ws.Range("J17").Value = 607.8
$ws.Range("L17").Value = 1823.4
$ws.Range("N17").Value = -2161.4

$ws.Range("H63").Value = 856
$ws.Range("I63").Value = 856
$ws.Range("K63").Value = 2568
$ws.Range("M63").Value = -1819

$ws.Range("H64").Value = 1450
$ws.Range("I64").Value = 1450
$ws.Range("K64").Value = 4350
$ws.Range("M64").Value = -4080

$ws.Range("H66").Value = 856
$ws.Range("I66").Value = 856
$ws.Range("K66").Value = 7704
$ws.Range("M66").Value = -3960

$ws.Range("H67").Value = 1450
$ws.Range("I67").Value = 1450
$ws.Range("K67").Value = 4350
$ws.Range("M67").Value = -3414

$ws.Range("H113").Value = 480.625
$ws.Range("I113").Value = 575
$ws.Range("J113").Value = 197.5
$ws.Range("K113").Value = 1725
$ws.Range("L113").Value = 592.5
$ws.Range("M113").Value = 445
$ws.Range("N113").Value = -4932.5

$ws.Range("H141").Value = 5919
$ws.Range("I141").Value = 7266.6665
$ws.Range("J141").Value = 3897.5
$ws.Range("K141").Value = 21799.9995
$ws.Range("L141").Value = 11692.5
$ws.Range("M141").Value = -16619.9995
$ws.Range("N141").Value = -22052.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 8666.666999999999
$ws.Range("I33").Value = 9400
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 9400
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = -9148
$ws.Range("N33").Value = -5504

$ws.Range("H97").Value = 713
$ws.Range("I97").Value = 277.25
$ws.Range("K97").Value = 277.25
$ws.Range("M97").Value = 218.75

$ws.Range("H117").Value = 24000
$ws.Range("J117").Value = 24000
$ws.Range("L117").Value = 24000
$ws.Range("N117").Value = -30884

$ws.Range("H133").Value = 120780
$ws.Range("J133").Value = 120780
$ws.Range("L133").Value = 120780
$ws.Range("N133").Value = -130900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15247.762
$ws.Range("I7").Value = 15233.777
$ws.Range("J7").Value = 15331.667
$ws.Range("K7").Value = 15233.777
$ws.Range("L7").Value = 15331.667
$ws.Range("M7").Value = -15121.777
$ws.Range("N7").Value = -15555.667

$ws.Range("H30").Value = 690.125
$ws.Range("I30").Value = 502.2857
$ws.Range("J30").Value = 2005
$ws.Range("K30").Value = 502.2857
$ws.Range("L30").Value = 2005
$ws.Range("M30").Value = -394.2857
$ws.Range("N30").Value = -2221

$ws.Range("H55").Value = 667.05884
$ws.Range("J55").Value = 1285
$ws.Range("L55").Value = 1285
$ws.Range("N55").Value = -1631

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws.Range("H126").Value = 15247.762
$ws.Range("I126").Value = 15233.777
$ws.Range("J126").Value = 15331.667
$ws.Range("K126").Value = 45701.331
$ws.Range("L126").Value = 45995.001
$ws.Range("M126").Value = -43231.331
$ws.Range("N126").Value = -50935.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 50007
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 50007
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 50007
$ws.Range("N61").Value = -50591
$ws.Range("M61").ClearContents()

$ws.Range("H69").Value = 11390.556
$ws.Range("J69").Value = 11390.556
$ws.Range("L69").Value = 11390.556
$ws.Range("N69").Value = -12888.556

$ws.Range("H72").Value = 11390.556
$ws.Range("J72").Value = 11390.556
$ws.Range("L72").Value = 34171.66800000001
$ws.Range("N72").Value = -41659.66800000001

$ws.Range("H126").Value = 4833.5625
$ws.Range("I126").Value = 3639.7273
$ws.Range("K126").Value = 10919.1819
$ws.Range("M126").Value = -8449.1819

Write-Host "Sheets updated via scheduled runner."
